# Generate Report for Handoff
#
# A new localization handoff completed for the f75f3b88-... file (row 7 of
# each per-locale sheet). Refresh the recorded "Latest Handoff Datetime"
# for both locales, and the rollup "Latest HO Xliff Generate Date" on the
# Overview sheet, to reflect the new handoff timestamps.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 7 is the
# f75f3b88-0e2f-4e5a-bc69-65700b6a22bc.md entry.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-17 08:41:28"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 7 is the same entry.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-17 08:41:22"

# de-de sheet: column H = "Latest Handoff Datetime", row 7 is the same entry.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-17 08:41:28"
